# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp text (row 1)
$ws.Range("A1").Value = "Datos actualizados a 11 de Abril de 2020 a las 12:22"

# Madrid row (row 4) - updated totals
$ws.Range("B4").Value = 45849
$ws.Range("C4").Value = 23663
$ws.Range("D4").Value = 16102
$ws.Range("E4").Value = 6084

# Navarra row (row 10) - updated totals
$ws.Range("B10").Value = 3817
$ws.Range("C10").Value = 603
$ws.Range("D10").Value = 2987
$ws.Range("E10").Value = 227

# La Rioja row (row 12) - updated totals
$ws.Range("B12").Value = 3223
$ws.Range("C12").Value = 1337
$ws.Range("D12").Value = 1679
$ws.Range("E12").Value = 207

# Reorder Cantabria to appear before Granada/Segovia, with fresh data,
# while Granada and Segovia keep their previous numbers (shifted down a row).
# First remove the old Cantabria row (row 27), shifting Segovia/Caceres/etc up.
$ws.Rows("27").Delete()

# Then insert a brand-new row at position 25 (before Granada) for the new Cantabria entry.
$ws.Rows("25").Insert()
$ws.Range("A25").Value = "Cantabria"
$ws.Range("B25").Value = 1719
$ws.Range("C25").Value = 281
$ws.Range("D25").Value = 1331
$ws.Range("E25").Value = 107
